# Updates the cryptos list price/volume figures (and two reordered rows)
# to match the latest GitHub Actions scrape.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column holds plain text (often dot-grouped numbers, e.g.
# "63.912.31") rather than real numbers. Forcing text format before the
# assignment -- then resetting the style -- stops Excel's automatic
# number detection from converting look-alike values (e.g. "589.16",
# "0.0000245") into numeric cells.
function Set-Text($RangeRef, $Text) {
    $cell = $ws.Range($RangeRef)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.Style = "Normal"
}

function Set-Row($Row, $D, $E) {
    if ($D -ne $null) {
        Set-Text "D$Row" $D
    }
    if ($E -ne $null) {
        $ws.Range("E$Row").Value = "  $E  "
    }
}

Set-Row 2  "63.912.31"  "+0.12%"
Set-Row 3  "3.137.46"   "+0.63%"
Set-Row 4  $null        "+0.10%"
Set-Row 5  "589.16"     "+0.45%"
Set-Row 6  "145.26"     "-0.63%"
Set-Row 7  $null        "+0.05%"
Set-Row 8  "3.129.81"   "+0.56%"
Set-Row 9  $null        "-0.28%"
Set-Row 10 $null        "-1.17%"
Set-Row 11 $null        "+2.52%"
Set-Row 12 $null        "-1.69%"
Set-Row 13 "0.0000245"  "-2.56%"
Set-Row 14 "37.27"      "+0.76%"
Set-Row 15 "3.656.91"   "+0.59%"
Set-Row 17 "7.31"       "+2.54%"
Set-Row 18 "63.765.12"  "+0.08%"
Set-Row 19 "3.137.64"   "+0.70%"
Set-Row 20 "467.97"     "+0.76%"
Set-Row 21 "14.33"      "+0.19%"
Set-Row 22 "0.731"      "+0.10%"
Set-Row 23 $null        "+0.07%"
Set-Row 24 "81.60"      "-0.57%"
Set-Row 25 $null        "-1.35%"
Set-Row 26 $null        "+6.94%"
Set-Row 27 $null        "+0.10%"
Set-Row 28 $null        "+9.94%"
Set-Row 29 "7.45"       "+8.48%"

# Rows 30 and 31 swap coin identity (ImmutableX <-> PancakeSwap)
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-Text "D30" "2.71"
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("B31").Value = "ImmutableX"
$ws.Range("C31").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-Text "D31" "2.24"
$ws.Range("E31").Value = "  +0.33%  "

Set-Row 32 $null        "+0.10%"
Set-Row 33 "27.69"      "+2.56%"
Set-Row 34 "0.109"      "+0.51%"
$subscript3 = [char]0x2083
$pepePrice = "0.0{0}0843" -f $subscript3
Set-Row 35 $pepePrice "-3.71%"
Set-Row 36 $null        "+1.49%"
Set-Row 37 $null        "+1.10%"
Set-Row 38 $null        "-2.69%"
Set-Row 39 $null        "-5.90%"

# Rows 40 and 41 swap coin identity (Cosmos <-> OKB)
$ws.Range("B40").Value = "OKB"
$ws.Range("C40").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-Text "D40" "51.25"
$ws.Range("E40").Value = "  +0.52%  "

$ws.Range("B41").Value = "Cosmos"
$ws.Range("C41").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Text "D41" "9.34"
$ws.Range("E41").Value = "  +7.49%  "

Set-Row 42 "454.30"     "+1.39%"
Set-Row 43 "0.293"      "+5.67%"
Set-Row 44 $null        "+0.09%"
Set-Row 45 "2.917.65"   "+1.08%"
Set-Row 46 "40.32"      "+12.05%"
Set-Row 47 $null        "-2.86%"
Set-Row 48 "132.73"     "+6.60%"
Set-Row 50 $null        "+2.53%"
Set-Row 51 $null        "-0.58%"
